$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$ws.Range("A5").Value = "Coconut"
$ws.Range("B5").Value = 22500

$ws.Range("A6").Value = "Black pepper"
$ws.Range("B6").Value = 36000
